# Applies the cryptos list update described by the commit diff.
# Strategy:
#  - Plain text cells (names, URLs, percentages) are written directly via .Value
#    since Excel will not reinterpret them as numbers.
#  - Cells whose new text is a bare numeric literal (e.g. "0.999", "8.25") would be
#    auto-converted to a Number by a plain .Value assignment, changing the stored cell
#    type away from the original inline string. To keep them as text (matching the
#    original t="inlineStr" cells) without adding any quote-prefix/text styling, we
#    write a formula that evaluates to the literal text, then convert that formula to
#    a static value in-place via Copy + PasteSpecial(xlPasteValues).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-NumericLookingText($addr, $val) {
    $ws.Range($addr).Formula = '="' + $val + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

Set-PlainText "D2" '61.937.04'
Set-PlainText "E2" '  +1.59%  '
Set-PlainText "D3" '2.405.91'
Set-PlainText "E3" '  +1.71%  '
Set-NumericLookingText "D4" '0.999'
Set-NumericLookingText "D5" '554.06'
Set-PlainText "E5" '  +1.32%  '
Set-NumericLookingText "D6" '142.19'
Set-PlainText "E6" '  +3.63%  '
Set-PlainText "E7" '  +0.08%  '
Set-NumericLookingText "D8" '0.530'
Set-PlainText "E8" '  +1.18%  '
Set-PlainText "D9" '2.398.40'
Set-PlainText "E9" '  +1.37%  '
Set-PlainText "E10" '  +1.61%  '
Set-PlainText "E11" '  -0.74%  '
Set-NumericLookingText "D12" '5.39'
Set-PlainText "E12" '  +1.01%  '
Set-NumericLookingText "D13" '0.351'
Set-PlainText "E13" '  +1.46%  '
Set-NumericLookingText "D14" '25.96'
Set-PlainText "E14" '  +4.26%  '
Set-NumericLookingText "D15" '0.0000174'
Set-PlainText "E15" '  +5.25%  '
Set-PlainText "D16" '2.853.82'
Set-PlainText "E16" '  +2.83%  '
Set-PlainText "D17" '61.960.48'
Set-PlainText "E17" '  +1.74%  '
Set-PlainText "D18" '2.409.54'
Set-PlainText "E18" '  +1.86%  '
Set-NumericLookingText "D19" '11.06'
Set-PlainText "E19" '  +3.17%  '
Set-NumericLookingText "D20" '4.18'
Set-PlainText "E20" '  +1.44%  '
Set-NumericLookingText "D21" '322.43'
Set-PlainText "E21" '  +0.90%  '
Set-NumericLookingText "D22" '6.69'
Set-PlainText "E22" '  +0.93%  '
Set-NumericLookingText "D23" '0.999'
Set-PlainText "E23" '  -0.09%  '
Set-NumericLookingText "D24" '64.96'
Set-PlainText "E24" '  +1.56%  '
Set-PlainText "E25" '  +5.02%  '
Set-NumericLookingText "D26" '8.97'
Set-PlainText "E26" '  +8.94%  '
Set-NumericLookingText "D27" '573.92'
Set-PlainText "E27" '  +15.69%  '
Set-PlainText "E28" '  +0.26%  '
Set-PlainText "D29" '2.525.07'
Set-PlainText "E29" '  +2.11%  '
Set-PlainText "B30" 'InternetComputer(DFINITY)'
Set-PlainText "C30" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-NumericLookingText "D30" '8.25'
Set-PlainText "E30" '  +2.06%  '
Set-PlainText "B31" 'PEPE'
Set-PlainText "C31" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-PlainText "D31" '0.0₃0924'
Set-PlainText "E31" '  +5.53%  '
Set-PlainText "B32" 'Fetch.AI'
Set-PlainText "C32" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-NumericLookingText "D32" '1.45'
Set-PlainText "E32" '  +5.77%  '
Set-NumericLookingText "D33" '0.148'
Set-PlainText "E33" '  -0.62%  '
Set-PlainText "E34" '  +2.33%  '
Set-NumericLookingText "D35" '1.56'
Set-PlainText "E35" '  +3.89%  '
Set-PlainText "E36" '  +0.09%  '
Set-NumericLookingText "D37" '5.63'
Set-PlainText "E37" '  +6.17%  '
Set-NumericLookingText "D38" '4.74'
Set-PlainText "E38" '  +1.73%  '
Set-PlainText "E39" '  +1.05%  '
Set-PlainText "B40" 'Monero'
Set-PlainText "C40" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-NumericLookingText "D40" '150.74'
Set-PlainText "E40" '  +3.89%  '
Set-PlainText "B41" 'EthereumClassic'
Set-PlainText "C41" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-NumericLookingText "D41" '18.58'
Set-PlainText "E41" '  +0.57%  '
Set-NumericLookingText "D42" '1.84'
Set-PlainText "E42" '  -1.97%  '
Set-PlainText "E43" '  +0.12%  '
Set-NumericLookingText "D44" '2.31'
Set-PlainText "E44" '  +13.93%  '
Set-NumericLookingText "D45" '148.89'
Set-PlainText "E45" '  +2.09%  '
Set-NumericLookingText "D46" '3.63'
Set-PlainText "E46" '  +1.67%  '
Set-NumericLookingText "D47" '0.0539'
Set-PlainText "E47" '  +3.66%  '
Set-PlainText "E48" '  +4.97%  '
Set-NumericLookingText "D49" '0.585'
Set-PlainText "E49" '  +2.28%  '
Set-NumericLookingText "D50" '0.0922'
Set-PlainText "E50" '  +1.87%  '
Set-PlainText "E51" '  +2.35%  '

$excel.CutCopyMode = 0

Write-Host "Applied cryptos list update."
